$d = $word.ActiveDocument

# 1. Title
$d.Paragraphs(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Musumeci attacca i suoi dipendenti: "L''80 per cento dei regionali si gratta la pancia"</w:t></w:r></w:p>')

# 2. Subtitle (add run)
$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>Il governatore torna a bacchettare gli impiegati: "Ora vogliono continuare il lavoro agile. Ma se non lavorate in ufficio come pensate di essere controllati a casa?". Il Siad-Cisal vuole quererarlo.</w:t></w:r></w:p>')

# 3. Big interview paragraph -> replaced
$d.Paragraphs(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:t>Nuova stoccata del presidente della Regione al personale dell''amministrazione da lui guidata. Intervenuto a Catania alle "Giornate dell''energia", il governatore siciliano Nello Musumeci ha preso di mira i dipendenti regionali, che definisce "improduttivi. L''80 per cento di loro si gratta la pancia dalla mattina alla sera". "Ma non ditelo ai sindacati - ha aggiunto - Ora vogliono stare ancora a casa per fare il ''lavoro agile''. Ma se non lavorate in ufficio, come pensate di essere controllati a casa?".</w:t><w:br/><w:br/><w:t xml:space="preserve">Giuseppe Badagliacca e Angelo Lo Curto del Siad-Cisal annunciano una querela nei confronti di Musumeci: “Abbiamo ascoltato con profondo sconcerto le sue parole contro i dipendenti regionali: accuse ingiuste, immotivate e offensive per tutti i lavoratori che ogni giorno svolgono il proprio dovere con abnegazione, anche in condizioni difficili. Evidentemente Musumeci è in difficoltà e prova a coprire i fallimenti del suo Governo puntando il dito contro l’anello più debole della catena, contro quei dipendenti che in piena pandemia hanno comunque lavorato e sono rientrati in servizio, nonostante la carenza dei dispositivi di sicurezza. Se la macchina non funziona non è colpa dei dipendenti, ma di chi politicamente ne è a capo. Adesso basta, la misura è colma: valuteremo con i nostri legali se sussistono gli estremi per una querela, tutelando i lavoratori in ogni sede”. </w:t><w:br/><w:t>try { MNZ_RICH(''Bottom''); } catch(e) {}</w:t></w:r></w:p>')

# 4. Word-definition entries (paragraphs 5..9)
$d.Paragraphs(5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>controllati</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>No pronunciation available</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>(verificare) (facts)</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')
$d.Paragraphs(6).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>ogni</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>[ˈoɲɲi]</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>(ciascuno)</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')
$d.Paragraphs(7).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>gratta</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>No pronunciation available</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>(sfregare)</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')
$d.Paragraphs(8).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>schifare</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>No pronunciation available</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>(disgustare, nauseare)</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')
$d.Paragraphs(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>magari</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>[maˈgari]</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>(forse, probabilmente)</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')

# 5. Append a brand new word-definition entry at the end (catania)
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br/><w:t xml:space="preserve">    word: </w:t><w:tab/><w:tab/><w:t>catania</w:t><w:br/><w:t xml:space="preserve">    pronunciation: </w:t><w:tab/><w:t>No pronunciation available</w:t><w:br/><w:t xml:space="preserve">    definition: </w:t><w:tab/><w:tab/><w:t>No additional information available</w:t><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p>')

Write-Output "edits applied"
